$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 0.9514656876067488

$ws.Range("D2").Value = $newValue
$ws.Range("D3").Value = $newValue
$ws.Range("D4").Value = $newValue
$ws.Range("D5").Value = $newValue
$ws.Range("D6").Value = $newValue
